# Automation Test Report - Build v4.22.00.159
# Update the "Employee" roster sheet: shift rows 6-8 down by one record,
# inserting the newest role-group/employee entry at row 6.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Employee")

$ws.Range("A6").Value = "33128455 - Mariela Gulgowski`nROLE : RTGO100 1701954977619"
$ws.Range("A7").Value = "ROLE GROUP : RTGO Operator 2023-12-07T19:27:58.156908600"
$ws.Range("A8").Value = "92970163 - Glenna Lynch`nROLE : RTGO100 1701853905917"
